# Update the "取得日時" (retrieved timestamp) column for the newly appended
# batch of rows (rows 2-12) on the "ランサーズ" sheet from the previous
# timestamp to the new append time: 2025-09-23 18:34:51

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-23 18:34:51"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
